# Applies the weekly fruit/vegetable data roll-forward for the
# "Mapocho Venta Directa de Santiago - Pepino ensalada" sheet.
# Row 2 takes the former Row 3 values, Row 3 takes the former Row 4
# values, and Row 4 takes the former Row 2 values for columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 3's data)
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("O2").Value = "Limache"
$ws.Range("P2").Value = 183

# Row 3 (was row 4's data)
$ws.Range("D3").Value = 44350
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 167

# Row 4 (was row 2's data)
$ws.Range("D4").Value = 44273
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 233
